# "Fruta / hortaliza, semanal" — weekly fruit/vegetable data refresh.
# Two new daily price records are inserted right after the header/lead-in
# block (before the existing row 438), pushing the rest of the table
# (old rows 438:533) down by two rows (to 440:535).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 438, shifting existing data down.
$ws.Range("A438:A439").EntireRow.Insert()

# --- New row 438 ---------------------------------------------------------
$ws.Range("A438").Value = 10
$ws.Range("B438").Value = "Vega Modelo de Temuco"
$ws.Range("C438").Value = "La Araucanía"
$ws.Range("D438").Value = 44943
$ws.Range("E438").Value = 9
$ws.Range("F438").Value = 100112024
$ws.Range("G438").Value = "Choclo"
$ws.Range("H438").Value = "Choclero"
$ws.Range("I438").Value = "Primera"
$ws.Range("J438").Value = 27000
$ws.Range("K438").Value = 280
$ws.Range("L438").Value = 300
$ws.Range("M438").Value = 291
$ws.Range("N438").Value = "$/unidad"
$ws.Range("O438").Value = "Región del Maule"
$ws.Range("P438").Value = 291
$ws.Range("Q438").Value = 1
$ws.Range("R438").Value = "Hortaliza"

# --- New row 439 ---------------------------------------------------------
$ws.Range("A439").Value = 10
$ws.Range("B439").Value = "Vega Modelo de Temuco"
$ws.Range("C439").Value = "La Araucanía"
$ws.Range("D439").Value = 44943
$ws.Range("E439").Value = 9
$ws.Range("F439").Value = 100112024
$ws.Range("G439").Value = "Choclo"
$ws.Range("H439").Value = "Dulce o Americano"
$ws.Range("I439").Value = "Primera"
$ws.Range("J439").Value = 15000
$ws.Range("K439").Value = 200
$ws.Range("L439").Value = 200
$ws.Range("M439").Value = 200
$ws.Range("N439").Value = "$/unidad"
$ws.Range("O439").Value = "Región del Maule"
$ws.Range("P439").Value = 200
$ws.Range("Q439").Value = 1
$ws.Range("R439").Value = "Hortaliza"
